$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @{
    "A2" = 1484.186284130758
    "B2" = 2238.001698579463
    "C2" = 2368.727307015083
    "D2" = 1899.401721890473
    "E2" = 1615.180614355501
    "F2" = 2354.115920891673
    "G2" = 672.4839635179353
    "H2" = 1693.022075552258
    "I2" = 2392.075819239393
    "J2" = 1427.478710753474
    "K2" = 1054.54181628627
    "L2" = 2183.684808525863
    "M2" = 951.1619631573842
    "N2" = 1912.987196290707
    "O2" = 614.4648568438683
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
